# Adds the "AdminSettings_User" test-data worksheet (with its table of
# user-creation test cases) to the workbook, placing it after the last
# existing sheet ("IMIA_Dashboard"), and nudges a couple of other sheets'
# selection state to mirror the authored workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the last tab in the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "AdminSettings_User"

# ---------------------------------------------------------------------
# 2. Header row (bold, like every other testdata sheet in the workbook).
# ---------------------------------------------------------------------
$headers = @("Description","TestType","FirstName","LastName","Email","Role","Success_or_Error?","Expected_Msg_FieldLevel","Expected_Msg_Header")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = [char](65 + $i)
    $ws.Range("$col" + "1").Value = $headers[$i]
}
$ws.Range("A1:I1").Font.Bold = $true

# ---------------------------------------------------------------------
# 3. Data rows (test cases for the Admin Settings "create user" screen).
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Select Role Empty"
$ws.Range("B2").Value = "Negative"
$ws.Range("C2").Value = "Ravi"
$ws.Range("D2").Value = "B"
$ws.Range("E2").Formula = '="Testinguser"&TEXT(NOW(),"ddmmhh")&"@gmail.com"'
$ws.Range("G2").Value = "Error"
$ws.Range("I2").Value = "Please select user's role"

$ws.Range("A3").Value = "FirstName is Empty"
$ws.Range("B3").Value = "Negative"
$ws.Range("D3").Value = "B"
$ws.Range("E3").Formula = '="Testinguser"&TEXT(NOW(),"ddmmhh")&"@gmail.com"'
$ws.Range("F3").Value = "User"
$ws.Range("G3").Value = "Error"
$ws.Range("I3").Value = "Please Enter First Name"

$ws.Range("A4").Value = "Special chars in FirstName"
$ws.Range("B4").Value = "Negative"
$ws.Range("C4").Value = "r@v8"
$ws.Range("D4").Value = "B"
$ws.Range("E4").Formula = '="Testinguser"&TEXT(NOW(),"ddmmhh")&"@gmail.com"'
$ws.Range("F4").Value = "User"
$ws.Range("G4").Value = "Error"
$ws.Range("I4").Value = "Please enter letters only"

$ws.Range("A5").Value = "Empty Email"
$ws.Range("B5").Value = "Negative"
$ws.Range("C5").Value = "Ravi"
$ws.Range("D5").Value = "B"
$ws.Range("F5").Value = "User"
$ws.Range("G5").Value = "Error"
$ws.Range("I5").Value = "Please Enter Email"

$ws.Range("A6").Value = "Invalid Email"
$ws.Range("B6").Value = "Negative"
$ws.Range("C6").Value = "Ravi"
$ws.Range("D6").Value = "B"
$ws.Range("E6").Value = "rav"
$ws.Range("F6").Value = "User"
$ws.Range("G6").Value = "Error"
$ws.Range("H6").Value = "Invalid Email"

$ws.Range("A7").Value = "Successful Creation_User"
$ws.Range("B7").Value = "Possitive"
$ws.Range("C7").Value = "Ravi"
$ws.Range("D7").Value = "B"
$ws.Range("E7").Formula = '="Testinguser0910111@gmail.com"'
$ws.Range("F7").Value = "User"
$ws.Range("G7").Value = "Success"
$ws.Range("I7").Value = "User created successfully. Email sent to user with credentials."

$ws.Range("A8").Value = "Successful Creation_SpAdmin"
$ws.Range("B8").Value = "Possitive"
$ws.Range("C8").Value = "Bolla"
$ws.Range("D8").Value = "R"
$ws.Range("E8").Formula = '="Testingspa"&TEXT(NOW(),"ddmmhh")&"@gmail.com"'
$ws.Range("F8").Value = "Specialty Admin"
$ws.Range("G8").Value = "Success"
$ws.Range("I8").Value = "User created successfully. Email sent to user with credentials."

# ---------------------------------------------------------------------
# 4. Column widths approximating the authored (auto-fit) widths.
# ---------------------------------------------------------------------
$colWidths = @{1=28.02; 2=17.02; 3=12.74; 4=13.02; 5=28.02; 6=15.88; 7=17.59; 8=23.59; 9=32.74}
foreach ($c in $colWidths.Keys) {
    $ws.Columns.Item($c).ColumnWidth = $colWidths[$c]
}

# ---------------------------------------------------------------------
# 5. Make this new sheet the active / visible tab, with the same
#    selected cell as in the authored workbook.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("E13").Select() | Out-Null

# ---------------------------------------------------------------------
# 6. Minor selection-state nudges on a couple of pre-existing sheets,
#    matching the authored workbook's last-saved cursor positions.
# ---------------------------------------------------------------------
$wsInvalidPassword = $wb.Worksheets.Item("Invalidpassword")
$wsInvalidPassword.Activate()
$wsInvalidPassword.Range("F1").Select() | Out-Null

$wsDashboard = $wb.Worksheets.Item("IMIA_Dashboard")
$wsDashboard.Activate()
$wsDashboard.Range("C1").Select() | Out-Null

# Leave the newly-added sheet as the active tab, as in the authored file.
$ws.Activate()
